$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "per 3 months" block to "per 4 months":
# Memory per 3 months (GB) -> Memory per 4 months (GB); formula I15/2 -> F15/3
$ws.Range("H19").Value = "Memory per 4 months (GB)"
$ws.Range("I19").Formula = "=F15/3"

# Battery capacity per 3 months (Ahr) -> Battery capacity per 4 months (Ahr); formula I16/2 -> F16/3
$ws.Range("H20").Value = "Battery capacity per 4 months (Ahr)"
$ws.Range("I20").Formula = "=F16/3"

# Update observed memory taken up (manual input) from 27 to 35
$ws.Range("I23").Value = 35

# Update the selected cell in the sheet view
$ws.Range("I24").Select()
